# Update "想去人数" (wish-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 2266
    "F3" = 1709
    "F5" = 1088
    "F6" = 811
    "F7" = 38
    "F8" = 5830
    "F9" = 88
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
